# Auto-generated Excel COM-interop script to apply the crypto list price/volume update.
# Source: diff of cryptos.xlsx (cell text values for columns B-E, rows 2-51).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Plain text updates (coin names, links, and already-non-numeric-looking price strings).
$plainUpdates = @{
    'D2' = '27.119.98'
    'E2' = '  -0.92%  '
    'D3' = '1.783.61'
    'E3' = '  -1.60%  '
    'E4' = '  -0.24%  '
    'E5' = '  -2.08%  '
    'E6' = '  -0.24%  '
    'E7' = '  +0.03%  '
    'E8' = '  -2.55%  '
    'E9' = '  -2.05%  '
    'E10' = '  -3.39%  '
    'E11' = '  -3.71%  '
    'E12' = '  -0.36%  '
    'E13' = '  -1.55%  '
    'E14' = '  -2.88%  '
    'D15' = '1.779.59'
    'E15' = '  -1.51%  '
    'E16' = '  -2.35%  '
    'E17' = '  -2.56%  '
    'E18' = '  -1.23%  '
    'E19' = '  -3.12%  '
    'E20' = '  -0.24%  '
    'E21' = '  -0.54%  '
    'E22' = '  -1.56%  '
    'D23' = '27.129.73'
    'E23' = '  -0.93%  '
    'E24' = '  -7.99%  '
    'E25' = '  -3.76%  '
    'E26' = '  -6.38%  '
    'E27' = '  -4.35%  '
    'E28' = '  -1.36%  '
    'E29' = '  +0.98%  '
    'D30' = '1.978.69'
    'E30' = '  -1.54%  '
    'E31' = '  -1.81%  '
    'E32' = '  -1.38%  '
    'E33' = '  -4.30%  '
    'E34' = '  -0.74%  '
    'E35' = '  -5.60%  '
    'E36' = '  -5.01%  '
    'E37' = '  -2.27%  '
    'E38' = '  -4.20%  '
    'E39' = '  -3.29%  '
    'E40' = '  -4.50%  '
    'E41' = '  -4.19%  '
    'B42' = 'TrustWalletToken'
    'C42' = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
    'E42' = '  -5.14%  '
    'B43' = 'FraxShare'
    'C43' = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
    'E43' = '  -6.61%  '
    'E44' = '  -3.60%  '
    'B45' = 'Decentraland'
    'C45' = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
    'E45' = '  -1.08%  '
    'B46' = 'Frax'
    'C46' = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
    'E46' = '  -0.31%  '
    'E47' = '  -4.22%  '
    'E48' = '  -2.65%  '
    'E49' = '  -0.52%  '
    'E50' = '  -3.33%  '
    'E51' = '  -2.36%  '
}

foreach ($ref in $plainUpdates.Keys) {
    $ws.Range($ref).Value = $plainUpdates[$ref]
}

# Price values in column D that look like plain numbers (e.g. '1.004') must be
# forced to Text format first, otherwise Excel auto-converts them to numeric
# values (losing the intended literal display text / introducing float drift).
$numericTextUpdates = @{
    'D4' = '1.004'
    'D5' = '336.76'
    'D6' = '1.001'
    'D7' = '0.3811'
    'D8' = '0.3412'
    'D9' = '48.15'
    'D10' = '1.189'
    'D11' = '0.07442'
    'D12' = '1.002'
    'D13' = '21.68'
    'D14' = '6.434'
    'D16' = '7.093'
    'D17' = '0.00001090'
    'D18' = '0.06644'
    'D19' = '83.33'
    'D20' = '1.001'
    'D21' = '6.525'
    'D22' = '17.38'
    'D24' = '12.23'
    'D25' = '2.372'
    'D26' = '2.500'
    'D27' = '21.12'
    'D28' = '1.458'
    'D29' = '154.97'
    'D31' = '133.89'
    'D32' = '3.991'
    'D33' = '6.051'
    'D34' = '0.08672'
    'D35' = '13.10'
    'D36' = '1.624'
    'D37' = '0.6840'
    'D38' = '5.391'
    'D39' = '0.06289'
    'D40' = '0.2177'
    'D41' = '0.02316'
    'D42' = '1.235'
    'D43' = '8.387'
    'D44' = '14.23'
    'D45' = '0.6441'
    'D46' = '0.9999'
    'D47' = '3.858'
    'D48' = '2.127'
    'D49' = '131.53'
    'D50' = '0.07098'
    'D51' = '78.60'
}

foreach ($ref in $numericTextUpdates.Keys) {
    $cell = $ws.Range($ref)
    $cell.NumberFormat = "@"
    $cell.Value = $numericTextUpdates[$ref]
}

